$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Reset baseline formatting for the BOM table (A1:F16) before rewriting
#    the contents, so every cell starts from a known style. Column A already
#    carries the "default" style (s=1, wrap text) used by most cells, and
#    E2 already carries the currency style (s=2) used by the Unit Cost /
#    Total Cost columns. Row 1 (headers) stays on the plain text style.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A1:F16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E2").Copy()
$ws.Range("E2:F16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Seed brand-new shared strings in the exact order they were first typed
#    by the author, so the shared string table matches the target order.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value  = "1276-1044-1-ND"
$ws.Range("B9").Value  = "1276-1000-1-ND"
$ws.Range("A8").Value  = "1uF ceramic"
$ws.Range("A10").Value = "10nF ceramic"
$ws.Range("B8").Value  = "1276-1866-1-ND"
$ws.Range("B10").Value = "732-8007-1-ND"
$ws.Range("B11").Value = "A130087CT-ND"
$ws.Range("B12").Value = "CR0603-JW-331ELFCT-ND"
$ws.Range("B13").Value = "RMCF0603JT10K0CT-ND"
$ws.Range("H11").Value = "Red 5mm THT LED"
$ws.Range("H12").Value = "Blue 5mm THT LED"
$ws.Range("H13").Value = "Yellow 5mm THT LED"
$ws.Range("H14").Value = "Green 5mm THT LED"
$ws.Range("F1").Value  = "Total Cost"

# ---------------------------------------------------------------------------
# 3. Rewrite the on-board parts table rows 1-16 (Part, Part Number, Need,
#    Have, Unit Cost). Total Cost (col F) is filled in afterwards with
#    formulas.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "On-board parts"
$ws.Range("B1").Value = "Part Number"
$ws.Range("C1").Value = "Needed"
$ws.Range("D1").Value = "Order"
$ws.Range("E1").Value = "Unit Cost"

$ws.Range("A2").Value = "ATMEGA"
$ws.Range("B2").Value = "ATMEGA32U4-AU"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 4.12

$ws.Range("A3").Value = "RFM69HCW - 915MHz"
$ws.Range("B3").Value = "1568-1394-ND"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 5.95

$ws.Range("A4").Value = "3.3V regulator"
$ws.Range("B4").Value = "296-39452-1-ND"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1.23

$ws.Range("A5").Value = "10uH inductor"
$ws.Range("B5").Value = "587-2886-1-ND"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 0.29

$ws.Range("A6").Value = "47uF ceramic"
$ws.Range("B6").Value = "587-1780-1-ND"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0.67

$ws.Range("A7").Value = "4.7uF ceramic"
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 0.15

$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = 0.034

$ws.Range("A9").Value = "0.1uF ceramic"
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 0.033

$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 0.019

$ws.Range("A11").Value = "220 ohm resistor"
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = 0.022

$ws.Range("A12").Value = "330 ohm resistor"
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = 0.015

$ws.Range("A13").Value = "10k resistor"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = 0.015

$ws.Range("A14").Value = "reset button"
$ws.Range("B14").Value = "401-1426-1-ND"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0.52

$ws.Range("A15").Value = "Level shifter"
$ws.Range("B15").Value = "296-12163-1-ND"
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 0.43

$ws.Range("A16").Value = "SMA PCB connector"
$ws.Range("B16").Value = "A97594-ND"
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 2.17

# ---------------------------------------------------------------------------
# 4. Total cost column: F2 stands alone, F3:F16 share one formula group,
#    and F18 sums the whole column.
# ---------------------------------------------------------------------------
$ws.Range("F2").Formula = "=D2*E2"
$ws.Range("F3:F16").Formula = "=D3*E3"

$ws.Range("E2").Copy()
$ws.Range("F18").PasteSpecial(-4122)   # xlPasteFormats -> s=2 (currency)
$excel.CutCopyMode = $false
$ws.Range("F18").Formula = "=SUM(F2:F16)"

# ---------------------------------------------------------------------------
# 5. Row 17 (old "50 ohm SMA coax" row) is removed entirely, and row 7's
#    special wrap style (Arial 7, used for long notes) moves to row 12
#    (the "Probably not using" style note moved off the SMA connector row).
# ---------------------------------------------------------------------------
$ws.Range("A17:M17").Clear()

$ws.Range("B7").Copy()
$ws.Range("B12").PasteSpecial(-4122)   # xlPasteFormats -> s=3 (Arial 7, wrap)
$excel.CutCopyMode = $false

$ws.Range("A1").Copy()
$ws.Range("B7").PasteSpecial(-4122)    # xlPasteFormats -> back to s=1
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6. Cosmetic sheet tweaks: column widths and row 15 height, selection.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 21.29
$ws.Columns.Item(6).ColumnWidth = 9.1
$ws.Rows.Item(15).AutoFit()

$ws.Range("F10").Select()

$excel.Calculate()
